# Refresh the crypto price snapshot (cryptos.xlsx) per the Mon Jul 22 06:45:26 UTC 2024
# GitHub Actions data-refresh commit: updates Price (D) / Volume(1h) (E) text cells,
# and for the two re-ranked coins (rows 25-26) also Coin (B) and Link (C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, never letting Excel auto-coerce
# numeric-looking strings (e.g. "0.588") into real numbers. We briefly flip the
# cell to the Text number format, assign the value, then restore the default
# "Normal" style so the cell ends up styled exactly as before the write.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '67.278.69'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '3.470.64'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-TextValue $ws.Range("D5") '593.16'
$ws.Range("E5").Value = '  -0.04%  '
Set-TextValue $ws.Range("D6") '179.08'
$ws.Range("E6").Value = '  +4.20%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '3.473.56'
$ws.Range("E8").Value = '  -0.41%  '
Set-TextValue $ws.Range("D9") '0.588'
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("E10").Value = '  +5.35%  '
Set-TextValue $ws.Range("D11") '7.08'
$ws.Range("E11").Value = '  -2.43%  '
Set-TextValue $ws.Range("D12") '0.431'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '4.074.19'
$ws.Range("E13").Value = '  -0.46%  '
Set-TextValue $ws.Range("D14") '32.13'
$ws.Range("E14").Value = '  +11.64%  '
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '67.278.58'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").Value = '3.481.54'
$ws.Range("E18").Value = '  -0.47%  '
Set-TextValue $ws.Range("D19") '6.24'
$ws.Range("E19").Value = '  -0.53%  '
Set-TextValue $ws.Range("D20") '14.29'
$ws.Range("E20").Value = '  +1.77%  '
Set-TextValue $ws.Range("D21") '390.12'
$ws.Range("E21").Value = '  -0.58%  '
Set-TextValue $ws.Range("D22") '7.86'
$ws.Range("E22").Value = '  -1.47%  '
Set-TextValue $ws.Range("D23") '72.82'
$ws.Range("E23").Value = '  +0.27%  '
Set-TextValue $ws.Range("D24") '0.997'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("B25").Value = 'LEO'
$ws.Range("C25").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D25") '5.71'
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("B26").Value = 'Polygon'
$ws.Range("C26").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D26") '0.534'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("E27").Value = '  +0.66%  '
Set-TextValue $ws.Range("D28") '10.35'
$ws.Range("E28").Value = '  +1.70%  '
Set-TextValue $ws.Range("D29") '0.175'
$ws.Range("E29").Value = '  -3.15%  '
$ws.Range("E30").Value = '  +0.49%  '
Set-TextValue $ws.Range("D31") '6.18'
$ws.Range("E31").Value = '  -0.64%  '
Set-TextValue $ws.Range("D32") '1.42'
$ws.Range("E32").Value = '  -0.28%  '
Set-TextValue $ws.Range("D33") '2.05'
$ws.Range("E33").Value = '  +0.12%  '
Set-TextValue $ws.Range("D34") '23.50'
$ws.Range("E34").Value = '  -0.64%  '
Set-TextValue $ws.Range("D35") '7.35'
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  -2.09%  '
Set-TextValue $ws.Range("D38") '163.52'
$ws.Range("E38").Value = '  +0.09%  '
Set-TextValue $ws.Range("D39") '0.868'
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("E41").Value = '  +7.01%  '
Set-TextValue $ws.Range("D42") '6.85'
$ws.Range("E42").Value = '  -1.16%  '
$ws.Range("D43").Value = '2.843.07'
$ws.Range("E43").Value = '  +1.52%  '
Set-TextValue $ws.Range("D44") '4.61'
$ws.Range("E44").Value = '  -1.12%  '
Set-TextValue $ws.Range("D45") '25.95'
$ws.Range("E45").Value = '  -0.82%  '
Set-TextValue $ws.Range("D46") '0.0723'
$ws.Range("E46").Value = '  -2.46%  '
Set-TextValue $ws.Range("D47") '26.48'
$ws.Range("E47").Value = '  -2.88%  '
Set-TextValue $ws.Range("D48") '41.91'
$ws.Range("E48").Value = '  -1.62%  '
Set-TextValue $ws.Range("D49") '0.0298'
$ws.Range("E49").Value = '  -1.10%  '
Set-TextValue $ws.Range("D50") '336.94'
$ws.Range("E50").Value = '  -0.04%  '
$ws.Range("E51").Value = '  -1.96%  '
